$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the two worker rows (16 and 17): JEILER MORALES MORALES now listed
# first (row 16) and RICARDO BELLO BLANCO second (row 17). Also update
# JEILER's Salario Basico (G column) to 1000000.

$ws.Range("C16").Value = "1050957528"
$ws.Range("D16").Value = "JEILER MORALES MORALES"
$ws.Range("F16").Value = 14536
$ws.Range("G16").Value = 1000000

$ws.Range("C17").Value = "20287343"
$ws.Range("D17").Value = "RICARDO BELLO BLANCO"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 908526
